# Updated soil type data
# Adds a new "pro_usda_soil_order" column (with full USDA soil order names)
# to the "profile" sheet, inserted before the existing "pro_soil_taxon"
# column, and adds the corresponding controlled-vocabulary list (the 12
# USDA soil orders) as a new column on the "controlled vocabulary" sheet,
# inserted before the existing "pro_soil_taxon_sys" column. A list data
# validation on the new profile column references the new vocabulary range.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "controlled vocabulary" sheet: insert new column E
#    (pro_usda_soil_order) before the existing pro_soil_taxon_sys column.
# ---------------------------------------------------------------
$cv = $wb.Worksheets.Item("controlled vocabulary")
$cv.Columns.Item(5).Insert()
$cv.Cells.Item(2, 5).Value = "pro_usda_soil_order"

$soilOrders = @("Alfisols","Andisols","Aridisols","Entisols","Gelisols","Histosols","Inceptisols","Mollisols","Oxisols","Spodosols","Ultisols","Vertisols")
for ($i = 0; $i -lt $soilOrders.Length; $i++) {
    $cv.Cells.Item(4 + $i, 5).Value = $soilOrders[$i]
}

# ---------------------------------------------------------------
# 2. "profile" sheet: insert new column N (pro_usda_soil_order) before
#    the existing pro_soil_taxon column, and populate with the full
#    soil-order names corresponding to the old abbreviated codes.
# ---------------------------------------------------------------
$pro = $wb.Worksheets.Item("profile")
$pro.Columns.Item(14).Insert()
$pro.Cells.Item(1, 14).Value = "pro_usda_soil_order"

$pro.Cells.Item(4, 14).Value = "Alfisols"
$pro.Cells.Item(5, 14).Value = "Alfisols"
$pro.Cells.Item(6, 14).Value = "Spodosols"
$pro.Cells.Item(7, 14).Value = "Alfisols"
$pro.Cells.Item(9, 14).Value = "Alfisols"

# List data validation on the new column, referencing the vocabulary range.
$validationRange = $pro.Range("N4:N1048576")
$validationRange.Validation.Add(3, 1, 1, "='controlled vocabulary'!`$E`$4:`$E`$15")
$validationRange.Validation.IgnoreBlank = $true
$validationRange.Validation.InCellDropdown = $true
